{"js": "// Insert the extra guidance sentence right after \"method.\" in section A4,\n// before the line break that ends that paragraph.\nconst body = context.document.body;\n\nconst results = body.search(\"method.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the anchor text \"method.\" in the document.');\n}\n\n// There is exactly one \"method.\" occurrence (end of part A4); use the first hit.\nconst target = results.items[0];\n\nconst insertion =\n  \" Keep in mind that this is not asking you to discuss how you\\u2019d implement your changes. \" +\n  \"You are only meant to tell the reader how you will let your team know to expect the changes in section A3.\";\n\ntarget.insertText(insertion, \"After\");\nawait context.sync();\n", "ps1": "# Insert the extra guidance sentence right after \"method.\" in section A4,\n# before the line break that ends that paragraph.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"method.\")\nif (-not $found) {\n    throw 'Could not find the anchor text \"method.\" in the document.'\n}\n\n$apostrophe = [char]0x2019\n$insertion = \" Keep in mind that this is not asking you to discuss how you\" + $apostrophe + \"d implement your changes. You are only meant to tell the reader how you will let your team know to expect the changes in section A3.\"\n\n$rng.InsertAfter($insertion)\n"}
